$wb = $excel.ActiveWorkbook

# The original "Sheet1" (sheetId 1) is replaced by a freshly added sheet
# named "ValidLogin" -- this is how Excel ends up minting sheetId 2 for it.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "ValidLogin"

# Login test data: header row + one valid-login row.
$newSheet.Range("A1").Value = "UserName"
$newSheet.Range("B1").Value = "Password"
$newSheet.Range("A2").Value = "admin"
$newSheet.Range("B2").Value = "manager"

# Drop the old sheet now that its data lives on the new one.
$wb.Worksheets("Sheet1").Delete() | Out-Null

# View state: zoom + selection on the new active sheet.
$excel.ActiveWindow.Zoom = 175
$newSheet.Range("B3").Select() | Out-Null
